$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '54.909.63'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  +8.12%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.436.29'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  +9.32%  '

$ws.Range('E4').Value = '  -0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '479.69'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +12.74%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '138.96'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +20.26%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.40%  '

$ws.Range('E8').Value = '  +11.29%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.461.03'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +10.12%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0964'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +15.36%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.46'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +6.80%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.323'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +10.89%  '

$ws.Range('E13').Value = '  +2.65%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.866.64'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  +8.89%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '55.023.45'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  +7.93%  '

$ws.Range('E16').Value = '  +12.28%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000134'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  +19.78%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.461.36'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  +8.39%  '

$ws.Range('E19').Value = '  +13.72%  '

$ws.Range('B20').Value = 'BitcoinCash'
$ws.Range('C20').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '313.69'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +8.51%  '

$ws.Range('B21').Value = 'Chainlink'
$ws.Range('C21').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.82'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +16.19%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.993'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +0.51%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.62'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +14.44%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '57.19'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +9.69%  '

$ws.Range('B25').Value = 'Binance-PegBSC-USD'
$ws.Range('C25').Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '1.00'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.17%  '

$ws.Range('B26').Value = 'Polygon'
$ws.Range('C26').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.403'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +13.33%  '

$ws.Range('E27').Value = '  +21.07%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.556.49'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +8.62%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '7.33'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +11.05%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0772'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +25.45%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.997'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -0.16%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '148.66'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +5.24%  '

$ws.Range('E33').Value = '  +10.42%  '

$ws.Range('E34').Value = '  +15.09%  '

$ws.Range('E35').Value = '  +13.44%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.12'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  +16.58%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.849'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +10.55%  '

$ws.Range('E38').Value = '  +10.44%  '

$ws.Range('B39').Value = 'OKB'
$ws.Range('C39').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '33.37'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +6.10%  '

$ws.Range('B40').Value = 'FirstDigitalUSD'
$ws.Range('C40').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.992'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  -0.63%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.602'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +10.75%  '

$ws.Range('E42').Value = '  +13.20%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.0542'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +11.91%  '

$ws.Range('E44').Value = '  +15.33%  '

$ws.Range('B45').Value = 'Bittensor'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '256.62'
$ws.Range('D45').ClearFormats()
$ws.Range('E45').Value = '  +36.48%  '

$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.66'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +21.69%  '

$ws.Range('B47').Value = 'WhiteBITCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '10.13'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.17%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0891'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +12.43%  '

$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.932.12'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +4.80%  '

$ws.Range('B50').Value = 'VeChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0221'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +13.04%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '17.02'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +12.71%  '
